# Update "Product" column values on the ProductUrls sheet from
# "JamesMockService" to "JamesTestService" (add example for API test).

$wb = $excel.ActiveWorkbook

$wsProductUrls = $wb.Worksheets.Item("ProductUrls")

$wsProductUrls.Range("B2").Value = "JamesTestService"
$wsProductUrls.Range("B3").Value = "JamesTestService"
$wsProductUrls.Range("B4").Value = "JamesTestService"

# Make ProductUrls the active sheet, with B4 selected (matches final
# author selection state captured in the saved file).
$wsProductUrls.Activate()
$wsProductUrls.Range("B4").Select()
